$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for every data row (2..390).
# All of these were bulk-updated from 45175 to 45177 in the source data.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 390 }

$ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3)).Value = 45177
